# "adding averages and more checks"
# - PERIOD TO EXPIRE (col H) drops by 8 days and LAST UPDATE (col I) moves
#   from 08-Sep-2025 to 16-Sep-2025 for every training row (3-39).
# - Header row (row 2) and the title (A1) switch to bold white text
#   (previously the title was bold/size-14/black and the header was just
#   bold/black on its blue fill).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# --- Title (A1) & header row (row 2): bold white font, default size -----
# (the title used to be bold/size-14/black; it now matches the header's
# bold/white/default-size look used on the blue-filled header row). The
# two worksheets share the same style table, so both the Training
# Dashboard and the Exam Dashboard titles/headers need the same update.
$headerRanges = @{
    "Training Dashboard" = "A2:K2"
    "Exam Dashboard"     = "A2:G2"
}

foreach ($sheetName in $headerRanges.Keys) {
    $sh = $wb.Worksheets.Item($sheetName)

    $title = $sh.Range("A1")
    $title.Font.Bold = $true
    $title.Font.Size = 11
    $title.Font.Color = 16777215   # RGB(255,255,255) -> white

    $header = $sh.Range($headerRanges[$sheetName])
    $header.Font.Bold = $true
    $header.Font.Size = 11
    $header.Font.Color = 16777215  # RGB(255,255,255) -> white
}

# --- PERIOD TO EXPIRE (H) / LAST UPDATE (I) updates ----------------------
$newPeriod = @{
    3=672; 4=674; 5=672; 6=674; 7=672; 8=673; 9=674; 10=673; 11=674; 12=684;
    13=674; 14=675; 15=675; 16=678; 17=678; 18=675; 19=675; 20=686; 21=308;
    22=310; 23=309; 24=323; 25=331; 26=330; 27=330; 28=328; 29=334; 30=334;
    31=338; 32=341; 33=341; 34=342; 35=342; 36=343; 37=343; 38=344; 39=344
}

# Scratch cell, formatted as Text, used so the "16-Sep-2025" string gets
# written (via copy / paste-values) as a literal text value instead of
# being auto-parsed into a date serial by Excel.
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"
$scratch.Value = "16-Sep-2025"
$scratch.Copy()

foreach ($row in 3..39) {
    $ws.Cells.Item($row, 8).Value = $newPeriod[$row]
    $ws.Range("I$row").PasteSpecial(-4163) | Out-Null   # xlPasteValues
}

$scratch.Clear()
$excel.CutCopyMode = $false
